$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Day 22 - "Hands-on tools: Nmap, Wireshark, basic scans") is now Done,
# so uncheck "In Progress?" (D) and check "Done?" (E).
$ws.Range("C9").Value = "Done"
$ws.Range("D9").Value = "☐"
$ws.Range("E9").Value = "☑"

# Row 10 (Day 22 - "Do beginner labs on TryHackMe / HackTheBox") moves to
# "In Progress" with the "In Progress?" checkbox checked.
$ws.Range("C10").Value = "In Progress"
$ws.Range("D10").Value = "☑"

# Update the selection to reflect where the user left off editing.
$ws.Range("B10").Select()
